$d = $word.ActiveDocument

# 1. Date of payment
$d.Content.Find.Execute("15 de Diciembre del 2020", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "20 de Julio del 2021", 2)

# 2. Agency code/name
$d.Content.Find.Execute("(274 ) SAN LUIS JILOTEPEQUE", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "( 1 ) ANTIGUA GUATEMALA", 2)

# 3. Entity name
$d.Content.Find.Execute("Frander Ivan Pérez Juárez ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "GLORIA ESTEPHANY MONROY DE LEON", 2)

# 4. "ref" -> "teléfono"
$d.Content.Find.Execute("ref", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "teléfono", 2)

# 5. Reference number
$d.Content.Find.Execute("T202111380", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " 78323191", 2)

# 6. Boleta number
$d.Content.Find.Execute("79400792 ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "3788000038", 2)

# 7. Certification issue date
$d.Content.Find.Execute("28 de Julio del 2021", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "30 de Julio del 2021", 2)
